# Lab 13 format update
# - Bold/size-14 "section title" cells (A1, A6, A21, A35) and the
#   regular "sub header" cells switch their font from Calibri to Arial
#   (and drop the Calibri theme "minor" scheme binding).
# - The plain data cells under each table (rows 8-17 cols B:F and
#   rows 24-33 cols B:C) get a new 10pt Arial font, vertically centered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Section-title cells: bold, 14pt, black -> keep everything, just
#    change the font family to Arial.
# ---------------------------------------------------------------
$titleCells = @("A1", "A6", "A21", "A35")
foreach ($addr in $titleCells) {
    $ws.Range($addr).Font.Name = "Arial"
}

# ---------------------------------------------------------------
# 2) Sub-header cells: regular, 11pt, black -> change font family to
#    Arial as well.
# ---------------------------------------------------------------
$subHeaderRanges = @("A2:B3", "A7:F7", "A22:C23", "A36:E37")
foreach ($addr in $subHeaderRanges) {
    $ws.Range($addr).Font.Name = "Arial"
}

# ---------------------------------------------------------------
# 3) Data rows: new 10pt Arial font, vertically centered, no fill /
#    border (same as default) -- applies to the "top clients" tables.
# ---------------------------------------------------------------
$dataRanges = @("B8:F17", "B24:C33")
foreach ($addr in $dataRanges) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Arial"
    $r.Font.Size = 10
    $r.Font.Color = 0
    $r.VerticalAlignment = -4108  # xlCenter
}
